$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the latest quarter (row 28: 2025Q2) metrics
$ws.Range("C28").Value = 299
$ws.Range("D28").Value = 30
$ws.Range("E28").Value = 269
$ws.Range("F28").Value = 4.672897196261682
